$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Sheet, $Address, $Text) {
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" '62.154.20'
Set-TextCell $ws "E2" '  +0.68%  '

# Row 3
Set-TextCell $ws "D3" '2.411.93'
Set-TextCell $ws "E3" '  +0.01%  '

# Row 4
Set-TextCell $ws "E4" '  -0.11%  '

# Row 5
Set-TextCell $ws "D5" '562.97'
Set-TextCell $ws "E5" '  +1.75%  '

# Row 6
Set-TextCell $ws "D6" '142.50'
Set-TextCell $ws "E6" '  +0.22%  '

# Row 7
Set-TextCell $ws "E7" '  +0.14%  '

# Row 8
Set-TextCell $ws "D8" '0.530'
Set-TextCell $ws "E8" '  +1.24%  '

# Row 9
Set-TextCell $ws "D9" '2.409.37'
Set-TextCell $ws "E9" '  -0.26%  '

# Row 10
Set-TextCell $ws "D10" '0.109'
Set-TextCell $ws "E10" '  +1.36%  '

# Row 11
Set-TextCell $ws "E11" '  -1.94%  '

# Row 12
Set-TextCell $ws "D12" '5.31'
Set-TextCell $ws "E12" '  -1.23%  '

# Row 13
Set-TextCell $ws "D13" '0.351'
Set-TextCell $ws "E13" '  -0.25%  '

# Row 14
Set-TextCell $ws "D14" '25.63'
Set-TextCell $ws "E14" '  -1.58%  '

# Row 15
Set-TextCell $ws "D15" '0.0000174'
Set-TextCell $ws "E15" '  -0.21%  '

# Row 16
Set-TextCell $ws "D16" '2.849.04'
Set-TextCell $ws "E16" '  +0.05%  '

# Row 17
Set-TextCell $ws "D17" '62.017.29'
Set-TextCell $ws "E17" '  +0.81%  '

# Row 18
Set-TextCell $ws "D18" '2.412.38'
Set-TextCell $ws "E18" '  -0.16%  '

# Row 19
Set-TextCell $ws "D19" '11.32'
Set-TextCell $ws "E19" '  +1.78%  '

# Row 20
Set-TextCell $ws "D20" '322.34'
Set-TextCell $ws "E20" '  -0.19%  '

# Row 21
Set-TextCell $ws "D21" '4.16'
Set-TextCell $ws "E21" '  -0.12%  '

# Row 22
Set-TextCell $ws "D22" '6.84'
Set-TextCell $ws "E22" '  +2.11%  '

# Row 23
Set-TextCell $ws "E23" '  -0.21%  '

# Row 24
Set-TextCell $ws "D24" '66.07'
Set-TextCell $ws "E24" '  +2.55%  '

# Row 25
Set-TextCell $ws "D25" '1.70'
Set-TextCell $ws "E25" '  -2.38%  '

# Row 26
Set-TextCell $ws "D26" '8.91'
Set-TextCell $ws "E26" '  -2.74%  '

# Row 27
Set-TextCell $ws "D27" '577.74'
Set-TextCell $ws "E27" '  +3.42%  '

# Row 28
Set-TextCell $ws "B28" 'WrappedeETH'
Set-TextCell $ws "C28" 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell $ws "D28" '2.534.54'
Set-TextCell $ws "E28" '  +1.35%  '

# Row 29
Set-TextCell $ws "B29" 'Binance-PegBSC-USD'
Set-TextCell $ws "C29" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws "D29" '1.00'
Set-TextCell $ws "E29" '  +0.25%  '

# Row 30
Set-TextCell $ws "B30" 'PEPE'
Set-TextCell $ws "C30" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws "D30" '0.0₃0943'
Set-TextCell $ws "E30" '  +2.50%  '

# Row 31
Set-TextCell $ws "D31" '8.21'
Set-TextCell $ws "E31" '  -0.92%  '

# Row 32
Set-TextCell $ws "D32" '1.42'
Set-TextCell $ws "E32" '  -1.06%  '

# Row 33
Set-TextCell $ws "E33" '  +0.20%  '

# Row 34
Set-TextCell $ws "E34" '  +1.02%  '

# Row 35
Set-TextCell $ws "E35" '  -0.34%  '

# Row 36
Set-TextCell $ws "E36" '  +0.20%  '

# Row 37
Set-TextCell $ws "B37" 'NEARProtocol'
Set-TextCell $ws "C37" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws "D37" '4.69'
Set-TextCell $ws "E37" '  -1.55%  '

# Row 38
Set-TextCell $ws "B38" 'RenderToken'
Set-TextCell $ws "C38" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell $ws "D38" '5.49'
Set-TextCell $ws "E38" '  -5.07%  '

# Row 39
Set-TextCell $ws "D39" '0.381'
Set-TextCell $ws "E39" '  -0.62%  '

# Row 40
Set-TextCell $ws "D40" '151.59'
Set-TextCell $ws "E40" '  +3.52%  '

# Row 41
Set-TextCell $ws "D41" '18.67'
Set-TextCell $ws "E41" '  -0.04%  '

# Row 42
Set-TextCell $ws "E42" '  -7.85%  '

# Row 43
Set-TextCell $ws "D43" '0.995'
Set-TextCell $ws "E43" '  -0.44%  '

# Row 44
Set-TextCell $ws "D44" '2.28'
Set-TextCell $ws "E44" '  +1.19%  '

# Row 45
Set-TextCell $ws "D45" '148.01'
Set-TextCell $ws "E45" '  -0.59%  '

# Row 46
Set-TextCell $ws "D46" '3.65'
Set-TextCell $ws "E46" '  +0.66%  '

# Row 47
Set-TextCell $ws "D47" '0.0533'
Set-TextCell $ws "E47" '  -0.28%  '

# Row 48
Set-TextCell $ws "D48" '19.93'
Set-TextCell $ws "E48" '  -1.65%  '

# Row 49
Set-TextCell $ws "D49" '0.593'
Set-TextCell $ws "E49" '  +0.36%  '

# Row 50
Set-TextCell $ws "D50" '0.0916'
Set-TextCell $ws "E50" '  +0.86%  '

# Row 51
Set-TextCell $ws "D51" '0.0227'
Set-TextCell $ws "E51" '  +0.96%  '

